$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two trailing rows (old rows 10 and 11), shrinking the table
# from 8 iterations (A1:D11) down to 6 iterations (A1:D9).
$ws.Range("A10:D11").Delete() | Out-Null

# Update the remaining data rows (2-9) with the recalculated
# secant-method values.
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = -75

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = -71.6576262890778

$ws.Range("B4").Value = 180.513140029589
$ws.Range("C4").Value = -44.059683603117
$ws.Range("D4").Value = 0.950142133705476

$ws.Range("B5").Value = 454.331183985633
$ws.Range("C5").Value = -16.9007869850754
$ws.Range("D5").Value = 0.602683798972299

$ws.Range("B6").Value = 624.726192982313
$ws.Range("C6").Value = -2.92090787921529
$ws.Range("D6").Value = 0.272751504436287

$ws.Range("B7").Value = 660.327940364287
$ws.Range("C7").Value = -0.170222338927957
$ws.Range("D7").Value = 0.0539152521129624

$ws.Range("B8").Value = 662.531104855656
$ws.Range("C8").Value = -0.0016467980762087
$ws.Range("D8").Value = 0.0033253751789487

$ws.Range("B9").Value = 662.552627355466
$ws.Range("C9").Value = -0.000000922459150842769
$ws.Range("D9").Value = 0.0000324842117008596
